$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the old "Total" row (old row 26), pushing it (and
# everything below) down by two rows.
$ws.Rows("26:27").Insert()

Write-Host "done"
